# Apply the VIN upload test fixture update:
#  - Change the VIN value in A2:A5 from "CCCKN3DD&E" to "XXXKN3DD&E"
#  - Keep VERSION / SYMBOL_2000 labels as-is (underlying shared-string reorder only)
#  - Move the active selection from D12 to A5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "XXXKN3DD&E"
$ws.Range("A3").Value = "XXXKN3DD&E"
$ws.Range("A4").Value = "XXXKN3DD&E"
$ws.Range("A5").Value = "XXXKN3DD&E"

$ws.Range("A5").Select()
